# Apply scheduled profit-recalculation updates to the Goblin Profits workbook.
# Values below were produced by the upstream market-data refresh job; each
# entry is the authoritative new value for the given cell on the given sheet.
$wb = $excel.ActiveWorkbook

$updates = @{}

$updates['ALC'] = @{
    'H2' = 166666850
    'I2' = 234.33333
    'J2' = 333333470
    'K2' = 234.33333
    'L2' = 333333470
    'M2' = -121.33333
    'N2' = -333333696
    'H17' = 2051.7058
    'J17' = 2546.818
    'L17' = 7640.454000000001
    'N17' = -7976.454000000001
    'H55' = 989.9091
    'I55' = 150.5
    'K55' = 150.5
    'M55' = 63.5
    'H69' = 5171
    'I69' = 4013
    'J69' = 5750
    'K69' = 12039
    'L69' = 17250
    'M69' = -11165
    'N69' = -18998
    'H72' = 5171
    'I72' = 4013
    'J72' = 5750
    'K72' = 36117
    'L72' = 51750
    'M72' = -31749
    'N72' = -60486
    'H131' = 7941105.5
    'J131' = 10399.143
    'L131' = 31197.429
    'N131' = -41277.429
    'H133' = 121982.8
    'J133' = 121982.8
    'L133' = 121982.8
    'N133' = -132102.8
    'H138' = 2271866.5
    'J138' = 2853861.2
    'L138' = 8561583.600000001
    'N138' = -8571863.600000001
    'H141' = 5270.7085
    'I141' = 2853.1765
    'K141' = 8559.529500000001
    'M141' = -3379.529500000001
}

$updates['ARM'] = @{
    'H32' = 152281.2
    'I32' = 152281.2
    'K32' = 152281.2
    'M32' = -151994.2
    'H37' = 38998.332
    'J37' = 49997.5
    'L37' = 49997.5
    'N37' = -50543.5
    'H63' = 6481.636
    'I63' = 2374.5
    'J63' = 8828.571
    'K63' = 2374.5
    'L63' = 8828.571
    'M63' = -1688.5
    'N63' = -10200.571
    'H66' = 6481.636
    'I66' = 2374.5
    'J66' = 8828.571
    'K66' = 11872.5
    'L66' = 44142.855
    'M66' = -8440.5
    'N66' = -51006.855
    'H74' = 1858.8948
    'I74' = 1569.1063
    'J74' = 3220.9
    'K74' = 1569.1063
    'L74' = 3220.9
    'M74' = -695.1062999999999
    'N74' = -4968.9
    'H77' = 1858.8948
    'I77' = 1569.1063
    'J77' = 3220.9
    'K77' = 7845.531499999999
    'L77' = 16104.5
    'M77' = -3477.531499999999
    'N77' = -24840.5
    'H97' = 1447
    'I97' = 815.8
    'J97' = 2499
    'K97' = 815.8
    'L97' = 2499
    'M97' = -319.8
    'N97' = -3491
}

$updates['BSM'] = @{
    'H20' = 2006.5172
    'I20' = 2031.2273
    'K20' = 2031.2273
    'M20' = -1784.2273
    'H86' = 19231946
    'I86' = 1206.5625
    'J86' = 50001132
    'K86' = 1206.5625
    'L86' = 50001132
    'M86' = -83.5625
    'N86' = -50003378
    'H89' = 19231946
    'I89' = 1206.5625
    'J89' = 50001132
    'K89' = 6032.8125
    'L89' = 250005660
    'M89' = -416.8125
    'N89' = -250016892
    'H134' = 597290.9
    'I134' = 2102.8538
    'K134' = 6308.5614
    'M134' = -3773.5614
}

$updates['CRP'] = @{
    'H12' = 1549.8572
    'I12' = 974.8333
    'J12' = 5000
    'K12' = 974.8333
    'L12' = 5000
    'M12' = -804.8333
    'N12' = -5340
    'H132' = 1860.3429
    'I132' = 1260.5518
    'J132' = 4759.3335
    'K132' = 3781.6554
    'L132' = 14278.0005
    'M132' = -1251.6554
    'N132' = -19338.0005
}

$updates['CUL'] = @{
    'H107' = 2275.2
    'J107' = 1852.6111
    'L107' = 5557.8333
    'N107' = -9397.8333
    'H121' = 898.7143
    'J121' = 999.5
    'L121' = 2998.5
    'N121' = -5618.5
    'H122' = 1050.0769
    'J122' = 1113.5555
    'L122' = 10021.9995
    'N122' = -14921.9995
}

$updates['GSM'] = @{
    'H2' = 97.611115
    'I2' = 86.09999999999999
    'J2' = 112
    'K2' = 86.09999999999999
    'L2' = 112
    'M2' = 26.90000000000001
    'N2' = -338
    'H59' = 13999.333
    'J59' = 13999.333
    'L59' = 13999.333
    'N59' = -15165.333
    'H96' = 12000
    'J96' = 12000
    'L96' = 12000
    'N96' = -17492
    'H97' = 849.1667
    'I97' = 799
    'K97' = 799
    'M97' = -303
    'H99' = 5805.8887
    'I99' = 5805.8887
    'K99' = 5805.8887
    'M99' = -3559.8887
    'H122' = 8447.111000000001
    'I122' = 8803.267
    'K122' = 26409.801
    'M122' = -23959.801
}

$updates['LTW'] = @{
    'H22' = 2208.7693
    'I22' = 1259.4
    'J22' = 2434.8096
    'K22' = 1259.4
    'L22' = 2434.8096
    'M22' = -964.4000000000001
    'N22' = -3024.8096
    'H27' = 2208.7693
    'I27' = 1259.4
    'J27' = 2434.8096
    'K27' = 1259.4
    'L27' = 2434.8096
    'M27' = -1152.4
    'N27' = -2648.8096
    'H93' = 3746.1875
    'I93' = 1858.3182
    'K93' = 1858.3182
    'M93' = -610.3181999999999
    'H132' = 3646.3438
    'I132' = 2632.4138
    'J132' = 13447.667
    'K132' = 7897.241399999999
    'L132' = 40343.001
    'M132' = -5367.241399999999
    'N132' = -45403.001
}

$updates['WVR'] = @{
    'H95' = 29058.6
    'J95' = 29058.6
    'L95' = 29058.6
    'N95' = -34550.6
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($ref in $cellUpdates.Keys) {
        $ws.Range($ref).Value = $cellUpdates[$ref]
    }
}

Write-Host 'Applied' ($updates.Values | ForEach-Object { $_.Count } | Measure-Object -Sum).Sum 'cell updates across' $updates.Keys.Count 'sheets'
